$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells holding Price/Volume figures are stored as text in this workbook.
# Force text formatting only on the specific cells being updated so the
# new numeric-looking / percentage-looking values are not auto-converted
# into Excel numbers, and untouched cells are left completely unaffected.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("E50").NumberFormat = "@"

$ws.Range("D2").Value = "256.63"
$ws.Range("E2").Value = "-1.02%"
$ws.Range("D3").Value = "27.21"
$ws.Range("E3").Value = "-2.86%"
$ws.Range("D4").Value = "4.558"
$ws.Range("E4").Value = "-12.60%"
$ws.Range("D5").Value = "0.05904"
$ws.Range("E5").Value = "-0.34%"
$ws.Range("D6").Value = "6.624"
$ws.Range("E6").Value = "-1.41%"
$ws.Range("D7").Value = "0.8601"
$ws.Range("E7").Value = "-1.58%"
$ws.Range("D8").Value = "0.9304"
$ws.Range("E8").Value = "-6.04%"
$ws.Range("D9").Value = "0.1407"
$ws.Range("E9").Value = "-0.91%"
$ws.Range("D10").Value = "0.03632"
$ws.Range("E10").Value = "-0.19%"
$ws.Range("D11").Value = "0.07089"
$ws.Range("E11").Value = "-2.09%"
$ws.Range("D12").Value = "0.03229"
$ws.Range("E12").Value = "0.83%"
$ws.Range("D13").Value = "0.09206"
$ws.Range("E13").Value = "-0.43%"
$ws.Range("D14").Value = "0.001559"
$ws.Range("E14").Value = "0.88%"
$ws.Range("D15").Value = "0.0006060"
$ws.Range("E15").Value = "-94.30%"
$ws.Range("D16").Value = "0.006086"
$ws.Range("E16").Value = "2.35%"
$ws.Range("D17").Value = "3.515"
$ws.Range("E17").Value = "0.52%"
$ws.Range("E18").Value = "-1.08%"
$ws.Range("E19").Value = "-0.14%"
$ws.Range("E20").Value = "-2.11%"
$ws.Range("E21").Value = "-0.95%"
$ws.Range("D22").Value = "3.850"
$ws.Range("E22").Value = "9.15%"
$ws.Range("D23").Value = "0.04220"
$ws.Range("E23").Value = "0.59%"
$ws.Range("D24").Value = "0.001218"
$ws.Range("E24").Value = "0.08%"
$ws.Range("D25").Value = "0.004279"
$ws.Range("E25").Value = "-6.43%"
$ws.Range("E26").Value = "0.17%"
$ws.Range("E27").Value = "0.07%"
$ws.Range("D40").Value = "0.03825"
$ws.Range("E40").Value = "-0.50%"
$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D41").Value = "0.006162"
$ws.Range("E41").Value = "13.85%"
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D42").Value = "0.1098"
$ws.Range("E42").Value = "-1.24%"
$ws.Range("D43").Value = "0.002201"
$ws.Range("E43").Value = "-7.27%"
$ws.Range("D44").Value = "0.01136"
$ws.Range("E44").Value = "3.99%"
$ws.Range("D45").Value = "0.00005435"
$ws.Range("E45").Value = "0.39%"
$ws.Range("E46").Value = "0.18%"
$ws.Range("D47").Value = "0.06020"
$ws.Range("E47").Value = "-29.50%"
$ws.Range("D48").Value = "0.08444"
$ws.Range("E48").Value = "3,849.14%"
$ws.Range("D49").Value = "0.00002101"
$ws.Range("E49").Value = "0.18%"
$ws.Range("D50").Value = "0.0002001"
$ws.Range("E50").Value = "0.18%"
